$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("data")
$wsPocetR = $wb.Worksheets.Item("pocetR")

# Fix header date typo (2021 -> 2022) on both sheets
$wsData.Range("W1").Value = "25. 1. 2022"
$wsPocetR.Range("V1").Value = "25. 1. 2022"

# Update "Nakaza v minulosti" percentage table on sheet "data" (rows 320-331)
# Row 320
$wsData.Range("D320").Value = 0
$wsData.Range("E320").Value = 0
$wsData.Range("F320").Value = 0
$wsData.Range("G320").Value = 0.03
$wsData.Range("H320").Value = 0.03
$wsData.Range("I320").Value = 0.04
$wsData.Range("J320").Value = 0.07000000000000001
$wsData.Range("K320").Value = 0.06
$wsData.Range("L320").Value = 0.09
$wsData.Range("M320").Value = 0.17
$wsData.Range("N320").Value = 0.36
$wsData.Range("O320").Value = 0.52
$wsData.Range("P320").Value = 0.61
$wsData.Range("Q320").Value = 0.68
$wsData.Range("R320").Value = 0.6899999999999999
$wsData.Range("S320").Value = 0.67
$wsData.Range("T320").Value = 0.67
$wsData.Range("U320").Value = 0.7
$wsData.Range("V320").Value = 0.72
$wsData.Range("W320").Value = 0.71

# Row 321
$wsData.Range("D321").Value = 0.27
$wsData.Range("E321").Value = 0.36
$wsData.Range("F321").Value = 0.34
$wsData.Range("G321").Value = 0.32
$wsData.Range("H321").Value = 0.33
$wsData.Range("I321").Value = 0.36
$wsData.Range("J321").Value = 0.31
$wsData.Range("K321").Value = 0.3
$wsData.Range("L321").Value = 0.25
$wsData.Range("M321").Value = 0.24
$wsData.Range("N321").Value = 0.17
$wsData.Range("O321").Value = 0.09
$wsData.Range("P321").Value = 0.03
$wsData.Range("Q321").Value = 0.02
$wsData.Range("R321").Value = 0.02
$wsData.Range("S321").Value = 0.02
$wsData.Range("T321").Value = 0.02
$wsData.Range("U321").Value = 0.02
$wsData.Range("V321").Value = 0
$wsData.Range("W321").Value = 0

# Row 322
$wsData.Range("D322").Value = 0.4
$wsData.Range("E322").Value = 0.24
$wsData.Range("F322").Value = 0.25
$wsData.Range("G322").Value = 0.27
$wsData.Range("H322").Value = 0.28
$wsData.Range("I322").Value = 0.29
$wsData.Range("J322").Value = 0.27
$wsData.Range("K322").Value = 0.3
$wsData.Range("L322").Value = 0.33
$wsData.Range("M322").Value = 0.26
$wsData.Range("N322").Value = 0.2
$wsData.Range("O322").Value = 0.1
$wsData.Range("P322").Value = 0.07000000000000001
$wsData.Range("Q322").Value = 0.08
$wsData.Range("R322").Value = 0.03
$wsData.Range("S322").Value = 0.04
$wsData.Range("T322").Value = 0.05
$wsData.Range("U322").Value = 0.03
$wsData.Range("V322").Value = 0.03
$wsData.Range("W322").Value = 0.04

# Row 323
$wsData.Range("D323").Value = 0
$wsData.Range("E323").Value = 0.07000000000000001
$wsData.Range("F323").Value = 0.08
$wsData.Range("G323").Value = 0.1
$wsData.Range("H323").Value = 0.12
$wsData.Range("I323").Value = 0.11
$wsData.Range("J323").Value = 0.13
$wsData.Range("K323").Value = 0.15
$wsData.Range("L323").Value = 0.14
$wsData.Range("M323").Value = 0.13
$wsData.Range("N323").Value = 0.13
$wsData.Range("O323").Value = 0.13
$wsData.Range("P323").Value = 0.11
$wsData.Range("Q323").Value = 0.09
$wsData.Range("R323").Value = 0.08
$wsData.Range("S323").Value = 0.08
$wsData.Range("T323").Value = 0.07000000000000001
$wsData.Range("U323").Value = 0.08
$wsData.Range("V323").Value = 0.06
$wsData.Range("W323").Value = 0.06

# Row 324
$wsData.Range("D324").Value = 0.16
$wsData.Range("E324").Value = 0.14
$wsData.Range("F324").Value = 0.21
$wsData.Range("G324").Value = 0.15
$wsData.Range("H324").Value = 0.15
$wsData.Range("I324").Value = 0.15
$wsData.Range("J324").Value = 0.13
$wsData.Range("K324").Value = 0.13
$wsData.Range("L324").Value = 0.12
$wsData.Range("M324").Value = 0.12
$wsData.Range("N324").Value = 0.09
$wsData.Range("O324").Value = 0.09
$wsData.Range("P324").Value = 0.1
$wsData.Range("Q324").Value = 0.06
$wsData.Range("R324").Value = 0.08
$wsData.Range("S324").Value = 0.1
$wsData.Range("T324").Value = 0.05
$wsData.Range("U324").Value = 0.09
$wsData.Range("V324").Value = 0.09
$wsData.Range("W324").Value = 0.1

# Row 325
$wsData.Range("D325").Value = 0.17
$wsData.Range("E325").Value = 0.19
$wsData.Range("F325").Value = 0.12
$wsData.Range("G325").Value = 0.13
$wsData.Range("H325").Value = 0.09
$wsData.Range("I325").Value = 0.05
$wsData.Range("J325").Value = 0.09
$wsData.Range("K325").Value = 0.06
$wsData.Range("L325").Value = 0.07000000000000001
$wsData.Range("M325").Value = 0.08
$wsData.Range("N325").Value = 0.05
$wsData.Range("O325").Value = 0.07000000000000001
$wsData.Range("P325").Value = 0.08
$wsData.Range("Q325").Value = 0.07000000000000001
$wsData.Range("R325").Value = 0.1
$wsData.Range("S325").Value = 0.09
$wsData.Range("T325").Value = 0.14
$wsData.Range("U325").Value = 0.08
$wsData.Range("V325").Value = 0.1
$wsData.Range("W325").Value = 0.09

# Row 326
$wsData.Range("D326").Value = 0
$wsData.Range("E326").Value = 0
$wsData.Range("F326").Value = 0
$wsData.Range("G326").Value = 0.01
$wsData.Range("H326").Value = 0.015
$wsData.Range("I326").Value = 0.02
$wsData.Range("J326").Value = 0.07000000000000001
$wsData.Range("K326").Value = 0.1
$wsData.Range("L326").Value = 0.14
$wsData.Range("M326").Value = 0.28
$wsData.Range("N326").Value = 0.45
$wsData.Range("O326").Value = 0.59
$wsData.Range("P326").Value = 0.64
$wsData.Range("Q326").Value = 0.68
$wsData.Range("R326").Value = 0.6899999999999999
$wsData.Range("S326").Value = 0.71
$wsData.Range("T326").Value = 0.73
$wsData.Range("U326").Value = 0.79
$wsData.Range("V326").Value = 0.78
$wsData.Range("W326").Value = 0.79

# Row 327
$wsData.Range("D327").Value = 0.33
$wsData.Range("E327").Value = 0.26
$wsData.Range("F327").Value = 0.27
$wsData.Range("G327").Value = 0.3
$wsData.Range("H327").Value = 0.28
$wsData.Range("I327").Value = 0.34
$wsData.Range("J327").Value = 0.31
$wsData.Range("K327").Value = 0.29
$wsData.Range("L327").Value = 0.28
$wsData.Range("M327").Value = 0.18
$wsData.Range("N327").Value = 0.11
$wsData.Range("O327").Value = 0.03
$wsData.Range("P327").Value = 0.02
$wsData.Range("Q327").Value = 0.01
$wsData.Range("R327").Value = 0
$wsData.Range("S327").Value = 0.005
$wsData.Range("T327").Value = 0.005
$wsData.Range("U327").Value = 0.005
$wsData.Range("V327").Value = 0
$wsData.Range("W327").Value = 0

# Row 328
$wsData.Range("D328").Value = 0.23
$wsData.Range("E328").Value = 0.23
$wsData.Range("F328").Value = 0.23
$wsData.Range("G328").Value = 0.23
$wsData.Range("H328").Value = 0.22
$wsData.Range("I328").Value = 0.21
$wsData.Range("J328").Value = 0.2
$wsData.Range("K328").Value = 0.21
$wsData.Range("L328").Value = 0.19
$wsData.Range("M328").Value = 0.17
$wsData.Range("N328").Value = 0.1
$wsData.Range("O328").Value = 0.06
$wsData.Range("P328").Value = 0.03
$wsData.Range("Q328").Value = 0.03
$wsData.Range("R328").Value = 0.03
$wsData.Range("S328").Value = 0.02
$wsData.Range("T328").Value = 0.03
$wsData.Range("U328").Value = 0.02
$wsData.Range("V328").Value = 0.02
$wsData.Range("W328").Value = 0.01

# Row 329
$wsData.Range("D329").Value = 0.11
$wsData.Range("E329").Value = 0.12
$wsData.Range("F329").Value = 0.14
$wsData.Range("G329").Value = 0.15
$wsData.Range("H329").Value = 0.16
$wsData.Range("I329").Value = 0.16
$wsData.Range("J329").Value = 0.16
$wsData.Range("K329").Value = 0.14
$wsData.Range("L329").Value = 0.14
$wsData.Range("M329").Value = 0.12
$wsData.Range("N329").Value = 0.11
$wsData.Range("O329").Value = 0.09
$wsData.Range("P329").Value = 0.1
$wsData.Range("Q329").Value = 0.08
$wsData.Range("R329").Value = 0.09
$wsData.Range("S329").Value = 0.08
$wsData.Range("T329").Value = 0.07000000000000001
$wsData.Range("U329").Value = 0.06
$wsData.Range("V329").Value = 0.05
$wsData.Range("W329").Value = 0.04

# Row 330
$wsData.Range("D330").Value = 0.12
$wsData.Range("E330").Value = 0.16
$wsData.Range("F330").Value = 0.2
$wsData.Range("G330").Value = 0.16
$wsData.Range("H330").Value = 0.16
$wsData.Range("I330").Value = 0.13
$wsData.Range("J330").Value = 0.12
$wsData.Range("K330").Value = 0.11
$wsData.Range("L330").Value = 0.11
$wsData.Range("M330").Value = 0.11
$wsData.Range("N330").Value = 0.1
$wsData.Range("O330").Value = 0.1
$wsData.Range("P330").Value = 0.09
$wsData.Range("Q330").Value = 0.08
$wsData.Range("R330").Value = 0.07000000000000001
$wsData.Range("S330").Value = 0.06
$wsData.Range("T330").Value = 0.06
$wsData.Range("U330").Value = 0.04
$wsData.Range("V330").Value = 0.05
$wsData.Range("W330").Value = 0.04

# Row 331
$wsData.Range("D331").Value = 0.21
$wsData.Range("E331").Value = 0.23
$wsData.Range("F331").Value = 0.16
$wsData.Range("G331").Value = 0.15
$wsData.Range("H331").Value = 0.165
$wsData.Range("I331").Value = 0.14
$wsData.Range("J331").Value = 0.14
$wsData.Range("K331").Value = 0.15
$wsData.Range("L331").Value = 0.14
$wsData.Range("M331").Value = 0.14
$wsData.Range("N331").Value = 0.13
$wsData.Range("O331").Value = 0.13
$wsData.Range("P331").Value = 0.12
$wsData.Range("Q331").Value = 0.12
$wsData.Range("R331").Value = 0.12
$wsData.Range("S331").Value = 0.125
$wsData.Range("T331").Value = 0.105
$wsData.Range("U331").Value = 0.08500000000000001
$wsData.Range("V331").Value = 0.1
$wsData.Range("W331").Value = 0.12

# Update counts table on sheet "pocetR" (rows 55-56)
# Row 55
$wsPocetR.Range("C55").Value = 15
$wsPocetR.Range("D55").Value = 110
$wsPocetR.Range("E55").Value = 155
$wsPocetR.Range("F55").Value = 175
$wsPocetR.Range("G55").Value = 220
$wsPocetR.Range("H55").Value = 241
$wsPocetR.Range("I55").Value = 280
$wsPocetR.Range("J55").Value = 284
$wsPocetR.Range("K55").Value = 288
$wsPocetR.Range("L55").Value = 292
$wsPocetR.Range("M55").Value = 293
$wsPocetR.Range("N55").Value = 273
$wsPocetR.Range("O55").Value = 251
$wsPocetR.Range("P55").Value = 265
$wsPocetR.Range("Q55").Value = 262
$wsPocetR.Range("R55").Value = 268
$wsPocetR.Range("S55").Value = 273
$wsPocetR.Range("T55").Value = 322
$wsPocetR.Range("U55").Value = 366
$wsPocetR.Range("V55").Value = 418

# Row 56
$wsPocetR.Range("C56").Value = 2152
$wsPocetR.Range("D56").Value = 2045
$wsPocetR.Range("E56").Value = 2031
$wsPocetR.Range("F56").Value = 1956
$wsPocetR.Range("G56").Value = 1900
$wsPocetR.Range("H56").Value = 1889
$wsPocetR.Range("I56").Value = 1821
$wsPocetR.Range("J56").Value = 1777
$wsPocetR.Range("K56").Value = 1771
$wsPocetR.Range("L56").Value = 1737
$wsPocetR.Range("M56").Value = 1682
$wsPocetR.Range("N56").Value = 1631
$wsPocetR.Range("O56").Value = 1531
$wsPocetR.Range("P56").Value = 1636
$wsPocetR.Range("Q56").Value = 1593
$wsPocetR.Range("R56").Value = 1568
$wsPocetR.Range("S56").Value = 1436
$wsPocetR.Range("T56").Value = 1468
$wsPocetR.Range("U56").Value = 1401
$wsPocetR.Range("V56").Value = 1397

Write-Host "Applied ZBP_07b_ockovani updates"
